# Refresh the cryptos price/volume table (GitHub Actions data refresh).
# Price (D) and Volume(1h) (E) columns are updated for most rows; rows 37/38
# (Aptos vs NEARProtocol) swap rank position, so Coin/Link/Price/Volume are
# all rewritten for those two rows.
#
# Numeric-looking Price strings (e.g. "1.00", "0.0000176") are written with a
# leading apostrophe so Excel stores them as text instead of coercing them to
# a Double (which would strip trailing zeros / flip to scientific notation),
# then the cell style is reset to "Normal" so no stray quote-prefix format is
# left applied to the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.589.16"
$ws.Range("E2").Value = "  +1.07%  "

$ws.Range("D3").Value = "3.449.73"
$ws.Range("E3").Value = "  +1.85%  "

$ws.Range("D4").Value = "'1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.03%  "

$ws.Range("D5").Value = "'580.27"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.23%  "

$ws.Range("D6").Value = "'150.15"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +8.91%  "

$ws.Range("D7").Value = "3.450.96"
$ws.Range("E7").Value = "  +1.90%  "

$ws.Range("E8").Value = "  +0.07%  "

$ws.Range("E9").Value = "  +0.92%  "

$ws.Range("D10").Value = "'7.80"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.28%  "

$ws.Range("D11").Value = "'0.127"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.54%  "

$ws.Range("E12").Value = "  +1.91%  "

$ws.Range("D13").Value = "4.040.71"
$ws.Range("E13").Value = "  +1.91%  "

$ws.Range("D14").Value = "'28.02"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +7.53%  "

$ws.Range("E15").Value = "  -0.36%  "

$ws.Range("D16").Value = "'0.0000176"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.63%  "

$ws.Range("D17").Value = "3.450.98"
$ws.Range("E17").Value = "  +1.88%  "

$ws.Range("D18").Value = "61.720.00"
$ws.Range("E18").Value = "  +1.01%  "

$ws.Range("D19").Value = "'6.31"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +8.77%  "

$ws.Range("D20").Value = "'14.40"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.93%  "

$ws.Range("D21").Value = "'9.50"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.57%  "

$ws.Range("D22").Value = "'389.68"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.68%  "

$ws.Range("D23").Value = "'0.568"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.47%  "

$ws.Range("D24").Value = "3.588.13"
$ws.Range("E24").Value = "  +1.66%  "

$ws.Range("D25").Value = "'72.87"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.51%  "

$ws.Range("E26").Value = "  +0.08%  "

$ws.Range("D27").Value = "'5.77"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.60%  "

$ws.Range("E28").Value = "  -0.46%  "

$ws.Range("E29").Value = "  +4.19%  "

$ws.Range("D30").Value = "'7.75"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +3.26%  "

$ws.Range("E31").Value = "  -12.64%  "

$ws.Range("D32").Value = "'0.999"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.13%  "

$ws.Range("E33").Value = "  +1.37%  "

$ws.Range("E34").Value = "  +1.06%  "

$ws.Range("E35").Value = "  -0.04%  "

$ws.Range("D36").Value = "'24.03"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.49%  "

$ws.Range("B37").Value = "Aptos"
$ws.Range("C37").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D37").Value = "'7.07"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.76%  "

$ws.Range("B38").Value = "NEARProtocol"
$ws.Range("C38").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D38").Value = "'5.25"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.69%  "

$ws.Range("E39").Value = "  +1.33%  "

$ws.Range("D40").Value = "'166.90"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.33%  "

$ws.Range("E41").Value = "  +3.98%  "

$ws.Range("D42").Value = "'27.02"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +11.12%  "

$ws.Range("D43").Value = "'0.794"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.94%  "

$ws.Range("E44").Value = "  +2.39%  "

$ws.Range("D45").Value = "'1.00"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.01%  "

$ws.Range("D46").Value = "'42.35"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.68%  "

$ws.Range("E47").Value = "  +0.54%  "

$ws.Range("D48").Value = "2.604.46"
$ws.Range("E48").Value = "  +5.67%  "

$ws.Range("E49").Value = "  -2.66%  "

$ws.Range("D50").Value = "'6.95"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.11%  "

$ws.Range("D51").Value = "'23.20"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.18%  "
